$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.319.46'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.879.16'
$ws.Range('E3').Value = '  -1.29%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.12'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.683'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.53'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +4.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.355'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '53.59'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0738'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0974'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.50'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.68%  '
$ws.Range('D14').Value = '2.149.58'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.764'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.91'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('D17').Value = '1.890.39'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = '35.329.67'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.79'
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Value = '0.0₃0820'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '243.97'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.84'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.95'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('E24').Value = '  +6.94%  '
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.20'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.48%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.78'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.51'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.26'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.127'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('D31').Value = '4.128.45'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.71'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.29%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.02'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0592'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.14'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.58%  '
$ws.Range('E37').Value = '  -0.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.839'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.94'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.69%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0729'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +11.43%  '
$ws.Range('E41').Value = '  +3.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0217'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '95.93'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.90%  '
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('D45').Value = '1.302.59'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('E47').Value = '  +5.80%  '
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.94'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.49%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.21'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.79%  '
